$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 261.22223
$ws.Range("I6").Value = 288.875
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 866.625
$ws.Range("L6").Value = 120
$ws.Range("M6").Value = -754.625
$ws.Range("N6").Value = -344
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H33").Value = 77693200
$ws.Range("I33").Value = 1111340.9
$ws.Range("J33").Value = 250002370
$ws.Range("K33").Value = 1111340.9
$ws.Range("L33").Value = 250002370
$ws.Range("M33").Value = -1111111.9
$ws.Range("N33").Value = -250002828
$ws.Range("H40").Value = 1718.8
$ws.Range("I40").Value = 1500.5
$ws.Range("K40").Value = 1500.5
$ws.Range("M40").Value = -1325.5
$ws.Range("H129").Value = 2227.8
$ws.Range("I129").Value = 1676.8
$ws.Range("J129").Value = 2778.8
$ws.Range("K129").Value = 5030.4
$ws.Range("L129").Value = 8336.400000000001
$ws.Range("M129").Value = -30.39999999999964
$ws.Range("N129").Value = -18336.4
$ws.Range("H132").Value = 2412.2856
$ws.Range("I132").Value = 2412.2856
$ws.Range("K132").Value = 7236.8568
$ws.Range("M132").Value = -4706.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2160.8845
$ws.Range("I2").Value = 1442.238
$ws.Range("K2").Value = 1442.238
$ws.Range("M2").Value = -1329.238
$ws.Range("H45").Value = 1777.2941
$ws.Range("I45").Value = 1453.125
$ws.Range("J45").Value = 2065.4443
$ws.Range("K45").Value = 1453.125
$ws.Range("L45").Value = 2065.4443
$ws.Range("M45").Value = -1076.125
$ws.Range("N45").Value = -2819.4443
$ws.Range("H94").Value = 76362.73
$ws.Range("J94").Value = 76362.73
$ws.Range("L94").Value = 76362.73
$ws.Range("N94").Value = -78164.73
$ws.Range("H102").Value = 33283.367
$ws.Range("I102").Value = 8260.6
$ws.Range("J102").Value = 127118.75
$ws.Range("K102").Value = 8260.6
$ws.Range("L102").Value = 127118.75
$ws.Range("M102").Value = -6638.6
$ws.Range("N102").Value = -130362.75
$ws.Range("H116").Value = 2160.8845
$ws.Range("I116").Value = 1442.238
$ws.Range("K116").Value = 1442.238
$ws.Range("M116").Value = 851.7619999999999
$ws.Range("H124").Value = 23862.334
$ws.Range("J124").Value = 23862.334
$ws.Range("L124").Value = 23862.334
$ws.Range("N124").Value = -33682.334
$ws.Range("H132").Value = 7392.9375
$ws.Range("I132").Value = 3897.7273
$ws.Range("K132").Value = 11693.1819
$ws.Range("M132").Value = -9163.1819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2160.8845
$ws.Range("I3").Value = 1442.238
$ws.Range("K3").Value = 1442.238
$ws.Range("M3").Value = -1328.238
$ws.Range("H22").Value = 215
$ws.Range("I22").Value = 215
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 215
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -42
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 928.0909
$ws.Range("I94").Value = 1149
$ws.Range("K94").Value = 1149
$ws.Range("M94").Value = -698

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H94").Value = 1584.3125
$ws.Range("I94").Value = 1846.8889
$ws.Range("J94").Value = 1246.7142
$ws.Range("K94").Value = 1846.8889
$ws.Range("L94").Value = 1246.7142
$ws.Range("M94").Value = -1395.8889
$ws.Range("N94").Value = -2148.7142
$ws.Range("H122").Value = 2152.2307
$ws.Range("I122").Value = 2032.75
$ws.Range("J122").Value = 2343.4
$ws.Range("K122").Value = 6098.25
$ws.Range("L122").Value = 7030.200000000001
$ws.Range("M122").Value = -3648.25
$ws.Range("N122").Value = -11930.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66667210
$ws.Range("I4").Value = 90909490
$ws.Range("K4").Value = 272728470
$ws.Range("M4").Value = -272728358
$ws.Range("H9").Value = 8999.799999999999
$ws.Range("J9").Value = 8999.799999999999
$ws.Range("L9").Value = 26999.4
$ws.Range("N9").Value = -27447.4
$ws.Range("H12").Value = 400.22223
$ws.Range("I12").Value = 105.28571
$ws.Range("J12").Value = 587.9091
$ws.Range("K12").Value = 315.85713
$ws.Range("L12").Value = 1763.7273
$ws.Range("M12").Value = -142.85713
$ws.Range("N12").Value = -2109.7273
$ws.Range("H21").Value = 3333
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 285.33334
$ws.Range("I23").Value = 71.666664
$ws.Range("J23").Value = 392.16666
$ws.Range("K23").Value = 214.999992
$ws.Range("L23").Value = 1176.49998
$ws.Range("M23").Value = 20.00000800000001
$ws.Range("N23").Value = -1646.49998
$ws.Range("H40").Value = 133
$ws.Range("I40").Value = 46.5
$ws.Range("K40").Value = 186
$ws.Range("M40").Value = -117
$ws.Range("H44").Value = 449
$ws.Range("I44").Value = 449
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1347
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -949
$ws.Range("N44").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 632
$ws.Range("I97").Value = 584
$ws.Range("J97").Value = 717.3333
$ws.Range("K97").Value = 584
$ws.Range("L97").Value = 717.3333
$ws.Range("M97").Value = -88
$ws.Range("N97").Value = -1709.3333
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H107").Value = 836.58826
$ws.Range("I107").Value = 731.9091
$ws.Range("J107").Value = 1028.5
$ws.Range("K107").Value = 731.9091
$ws.Range("L107").Value = 1028.5
$ws.Range("M107").Value = 1188.0909
$ws.Range("N107").Value = -4868.5
$ws.Range("H132").Value = 5197.1
$ws.Range("I132").Value = 5197.1
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15591.3
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13061.3
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 127500
$ws.Range("J136").Value = 127500
$ws.Range("L136").Value = 382500
$ws.Range("N136").Value = -387600

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 8063.6665
$ws.Range("I40").Value = 8796.714
$ws.Range("J40").Value = 7037.4
$ws.Range("K40").Value = 8796.714
$ws.Range("L40").Value = 7037.4
$ws.Range("M40").Value = -8660.714
$ws.Range("N40").Value = -7309.4
$ws.Range("H55").Value = 1062.3
$ws.Range("I55").Value = 287.83334
$ws.Range("J55").Value = 2224
$ws.Range("K55").Value = 287.83334
$ws.Range("L55").Value = 2224
$ws.Range("M55").Value = -114.83334
$ws.Range("N55").Value = -2570
$ws.Range("H61").Value = 3639.875
$ws.Range("J61").Value = 2696.75
$ws.Range("L61").Value = 2696.75
$ws.Range("N61").Value = -3100.75
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H113").Value = 3639.875
$ws.Range("J113").Value = 2696.75
$ws.Range("L113").Value = 2696.75
$ws.Range("N113").Value = -7036.75
$ws.Range("H122").Value = 5964.0625
$ws.Range("I122").Value = 4830.5
$ws.Range("J122").Value = 7853.3335
$ws.Range("K122").Value = 14491.5
$ws.Range("L122").Value = 23560.0005
$ws.Range("M122").Value = -12041.5
$ws.Range("N122").Value = -28460.0005
$ws.Range("H136").Value = 2812.5
$ws.Range("J136").Value = 2875
$ws.Range("L136").Value = 8625
$ws.Range("N136").Value = -13725

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5698.7334
$ws.Range("I122").Value = 5837.5557
$ws.Range("J122").Value = 5490.5
$ws.Range("K122").Value = 17512.6671
$ws.Range("L122").Value = 16471.5
$ws.Range("M122").Value = -15062.6671
$ws.Range("N122").Value = -21371.5
